# Update the date header and the 25 division problems in the table to the
# "next day" worksheet values, per commit "Update master to output generated
# at c8c62b6".
$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-12-25 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-12-26 Friday", 2) | Out-Null
$d.Content.Find.Execute("40÷3=13, 1", $true, $false, $false, $false, $false, $true, 1, $false, "60÷4=15, 0", 2) | Out-Null
$d.Content.Find.Execute("96÷6=16, 0", $true, $false, $false, $false, $false, $true, 1, $false, "26÷2=13, 0", 2) | Out-Null
$d.Content.Find.Execute("14÷9=1, 5", $true, $false, $false, $false, $false, $true, 1, $false, "15÷7=2, 1", 2) | Out-Null
$d.Content.Find.Execute("38÷7=5, 3", $true, $false, $false, $false, $false, $true, 1, $false, "68÷5=13, 3", 2) | Out-Null
$d.Content.Find.Execute("81÷2=40, 1", $true, $false, $false, $false, $false, $true, 1, $false, "50÷4=12, 2", 2) | Out-Null
$d.Content.Find.Execute("20÷4=5, 0", $true, $false, $false, $false, $false, $true, 1, $false, "72÷5=14, 2", 2) | Out-Null
$d.Content.Find.Execute("95÷4=23, 3", $true, $false, $false, $false, $false, $true, 1, $false, "53÷3=17, 2", 2) | Out-Null
$d.Content.Find.Execute("40÷2=20, 0", $true, $false, $false, $false, $false, $true, 1, $false, "79÷5=15, 4", 2) | Out-Null
$d.Content.Find.Execute("78÷3=26, 0", $true, $false, $false, $false, $false, $true, 1, $false, "15÷4=3, 3", 2) | Out-Null
$d.Content.Find.Execute("54÷9=6, 0", $true, $false, $false, $false, $false, $true, 1, $false, "44÷5=8, 4", 2) | Out-Null
$d.Content.Find.Execute("84÷9=9, 3", $true, $false, $false, $false, $false, $true, 1, $false, "32÷2=16, 0", 2) | Out-Null
$d.Content.Find.Execute("59÷2=29, 1", $true, $false, $false, $false, $false, $true, 1, $false, "96÷3=32, 0", 2) | Out-Null
$d.Content.Find.Execute("23÷8=2, 7", $true, $false, $false, $false, $false, $true, 1, $false, "28÷5=5, 3", 2) | Out-Null
$d.Content.Find.Execute("71÷4=17, 3", $true, $false, $false, $false, $false, $true, 1, $false, "88÷2=44, 0", 2) | Out-Null
$d.Content.Find.Execute("96÷9=10, 6", $true, $false, $false, $false, $false, $true, 1, $false, "38÷3=12, 2", 2) | Out-Null
$d.Content.Find.Execute("90÷9=10, 0", $true, $false, $false, $false, $false, $true, 1, $false, "41÷7=5, 6", 2) | Out-Null
$d.Content.Find.Execute("77÷4=19, 1", $true, $false, $false, $false, $false, $true, 1, $false, "78÷7=11, 1", 2) | Out-Null
$d.Content.Find.Execute("81÷7=11, 4", $true, $false, $false, $false, $false, $true, 1, $false, "77÷7=11, 0", 2) | Out-Null
$d.Content.Find.Execute("46÷4=11, 2", $true, $false, $false, $false, $false, $true, 1, $false, "18÷8=2, 2", 2) | Out-Null
$d.Content.Find.Execute("85÷8=10, 5", $true, $false, $false, $false, $false, $true, 1, $false, "32÷4=8, 0", 2) | Out-Null
$d.Content.Find.Execute("98÷8=12, 2", $true, $false, $false, $false, $false, $true, 1, $false, "98÷2=49, 0", 2) | Out-Null
$d.Content.Find.Execute("11÷3=3, 2", $true, $false, $false, $false, $false, $true, 1, $false, "68÷7=9, 5", 2) | Out-Null
$d.Content.Find.Execute("77÷6=12, 5", $true, $false, $false, $false, $false, $true, 1, $false, "72÷9=8, 0", 2) | Out-Null
$d.Content.Find.Execute("65÷5=13, 0", $true, $false, $false, $false, $false, $true, 1, $false, "24÷6=4, 0", 2) | Out-Null
$d.Content.Find.Execute("10÷3=3, 1", $true, $false, $false, $false, $false, $true, 1, $false, "33÷3=11, 0", 2) | Out-Null
